$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewSavingInput")
$ws.Range("B5").Value = "save"
$ws.Activate()
$ws.Range("B13").Select()
